$d = $word.ActiveDocument

# --- 1. Student login JSON: add type, drop "s" prefix from sID/sPassword ---
# paragraph 2 (COM, 1-based) == "{<LDQ>sID<RDQ>:<LDQ>x<RDQ>, <LDQ>sPassword<RDQ>:<LDQ>x<RDQ>}    <CJK>student</CJK>"
$p1 = $d.Paragraphs(2).Range
$p1.Find.Execute([string]([char]0x201C) + "sID" + [char]0x201D, $true, $false, $false, $false, $false, $true, 1, $false, ([char]0x201C) + "type" + ([char]0x201D) + ":3," + ([char]0x201C) + "ID" + ([char]0x201D), 1)
$p1 = $d.Paragraphs(2).Range
$p1.Find.Execute([string]([char]0x201C) + "sPassword" + [char]0x201D, $true, $false, $false, $false, $false, $true, 1, $false, ([char]0x201C) + "Password" + ([char]0x201D), 1)

# --- 2. Teacher login JSON: add type, drop "t" prefix from tID/tPassword ---
$p2 = $d.Paragraphs(5).Range
$p2.Find.Execute([string]([char]0x201C) + "tID" + [char]0x201D, $true, $false, $false, $false, $false, $true, 1, $false, ([char]0x201C) + "type" + ([char]0x201D) + ":2," + ([char]0x201C) + "ID" + ([char]0x201D), 1)
$p2 = $d.Paragraphs(5).Range
$p2.Find.Execute([string]([char]0x201C) + "tPassword" + [char]0x201D, $true, $false, $false, $false, $false, $true, 1, $false, ([char]0x201C) + "Password" + ([char]0x201D), 1)

# --- 3. Admin-login paragraph: drop the JSON entirely, keep only "没做"; ---
#        merge away the now-duplicate "没做" paragraph + the blank line after it.
$p3 = $d.Paragraphs(8).Range
$p3.End = $p3.End - 1
$p3.Text = [string]([char]0x6CA1) + [string]([char]0x505A)

$p4 = $d.Paragraphs(9).Range
$p4.End = $p4.End - 1
$p4.Text = ""

$d.Paragraphs(9).Range.Delete()

# --- 4. "_GoBack" (last-edit marker) follows the edit: drop it from its old
#        spot after "/correct" and re-add it right after the new "没做" text ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$pEdit = $d.Paragraphs(8).Range
$goBackPos = $pEdit.End - 1
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
